# "Generate Report for Handoff"
# Replace the previous handback-report row data with a freshly generated
# handoff report: new source file (23a96089-...), a new second file
# (ffffe4f65d2d-...), status flips from "Handed back: in sync with en-US"
# to "Ready for handoff", timestamps refreshed, and the now-obsolete
# "Latest Target File" / "Dependency From" columns (F/G) are dropped from
# the per-language detail sheets.

$wb = $excel.ActiveWorkbook

$srcMd1 = "23a96089-9409-4a50-8a86-bca15a42054d.md"
$srcMd2 = "ffffe4f65d2d-4d6b-4114-b747-55508ca8ef7b.md"
$status = "Ready for handoff"
$overviewDate = "2016-50-17 16:50:00"

$xlfZh = "23a96089-9409-4a50-8a86-bca15a42054d.a68b0890a10c6cd7e927b1996654fa84ef4ea170.zh-cn.xlf"
$xlfDe = "23a96089-9409-4a50-8a86-bca15a42054d.a68b0890a10c6cd7e927b1996654fa84ef4ea170.de-de.xlf"
$handoffDtZh = "2016-03-17 16:49:57"
$handoffDtDe = "2016-03-17 16:50:00"
$handbackDt = "0001-01-01 00:00:00"

$baseCommit = "152df1c26c7fb1f343c88ceee4548f294e6f57c3"
$hoCommitZh = "ea7c097be589768129f45fd42c149c7f0d4ebe4c"
$hoCommitDe = "37ceb674b7db32d5f1ce8f30a2f50f8ce1ea21de"

$mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/$baseCommit/e2e/$srcMd1"
$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/$baseCommit/e2e/$srcMd2"
$xlfUrlZh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hoCommitZh/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZh"
$xlfUrlDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hoCommitDe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDe"

function Restyle-Hyperlink($range) {
    # Approximate the workbook's custom "HyperLink" look (underlined,
    # #6495ED) on cells touched by Hyperlinks.Add, which otherwise stamps
    # the generic theme hyperlink format.
    $range.Font.Underline = $true
    $range.Font.Color = 15570276
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Cells.Hyperlinks.Delete()

$ws.Range("A2").Value2 = $srcMd1
$ws.Range("B2").Value2 = $status
$ws.Range("C2").Value2 = $status
$ws.Range("D2").Value2 = $overviewDate

$ws.Range("A3").Value2 = $srcMd2
$ws.Range("B3").Value2 = $status
$ws.Range("C3").Value2 = $status
$ws.Range("D3").Value2 = $overviewDate

$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $srcMd1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $srcMd2) | Out-Null
Restyle-Hyperlink $ws.Range("A2")
Restyle-Hyperlink $ws.Range("A3")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Cells.Hyperlinks.Delete()
$ws.Range("F2:G3").Clear()

$ws.Range("A2").Value2 = $srcMd1
$ws.Range("B2").Value2 = ".md"
$ws.Range("C2").Value2 = $status
$ws.Range("D2").Value2 = $xlfZh
$ws.Range("E2").Value2 = $handoffDtZh
$ws.Range("H2").Value2 = $handbackDt

$ws.Range("A3").Value2 = $srcMd2
$ws.Range("B3").Value2 = ".md"
$ws.Range("C3").Value2 = $status
$ws.Range("D3").Value2 = $xlfZh
$ws.Range("E3").Value2 = $handoffDtZh
$ws.Range("H3").Value2 = $handbackDt

$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $srcMd1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl1, "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), $xlfUrlZh, "", "", $xlfZh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $srcMd2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl2, "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), $xlfUrlZh, "", "", $xlfZh) | Out-Null
Restyle-Hyperlink $ws.Range("A2")
Restyle-Hyperlink $ws.Range("B2")
Restyle-Hyperlink $ws.Range("D2")
Restyle-Hyperlink $ws.Range("A3")
Restyle-Hyperlink $ws.Range("B3")
Restyle-Hyperlink $ws.Range("D3")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Cells.Hyperlinks.Delete()
$ws.Range("F2:G3").Clear()

$ws.Range("A2").Value2 = $srcMd1
$ws.Range("B2").Value2 = ".md"
$ws.Range("C2").Value2 = $status
$ws.Range("D2").Value2 = $xlfDe
$ws.Range("E2").Value2 = $handoffDtDe
$ws.Range("H2").Value2 = $handbackDt

$ws.Range("A3").Value2 = $srcMd2
$ws.Range("B3").Value2 = ".md"
$ws.Range("C3").Value2 = $status
$ws.Range("D3").Value2 = $xlfDe
$ws.Range("E3").Value2 = $handoffDtDe
$ws.Range("H3").Value2 = $handbackDt

$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $srcMd1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl1, "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), $xlfUrlDe, "", "", $xlfDe) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $srcMd2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl2, "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), $xlfUrlDe, "", "", $xlfDe) | Out-Null
Restyle-Hyperlink $ws.Range("A2")
Restyle-Hyperlink $ws.Range("B2")
Restyle-Hyperlink $ws.Range("D2")
Restyle-Hyperlink $ws.Range("A3")
Restyle-Hyperlink $ws.Range("B3")
Restyle-Hyperlink $ws.Range("D3")
